$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-07-31"

# Update the "July (through 07-30)" label to "July (through 07-31)".
$ws.Range("A8").Value = "July (through 07-31)"

# Update the July row (row 8) values for 2016-2022.
$ws.Range("C8").Value = 54
$ws.Range("D8").Value = 75
$ws.Range("E8").Value = 72
$ws.Range("F8").Value = 53
$ws.Range("G8").Value = 149
$ws.Range("H8").Value = 150
$ws.Range("I8").Value = 169

# Update the Total row (row 9) values for 2016-2022.
$ws.Range("C9").Value = 302
$ws.Range("D9").Value = 465
$ws.Range("E9").Value = 425
$ws.Range("F9").Value = 304
$ws.Range("G9").Value = 621
$ws.Range("H9").Value = 910
$ws.Range("I9").Value = 975
